# Update root path for qualitron
# Qualitron_Ruta.xlsx - the "Red" column values in the last two data rows
# ("No") are corrected to "Si", and the active selection moves to B16.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 12 and 13 in column A ("Red") were "No" - change them to "Si".
$ws.Range("A12").Value = "Si"
$ws.Range("A13").Value = "Si"

# Update the active selection to match the saved cursor position.
$ws.Range("B16").Select()
